$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16: change the GUID of the table style applied to the small table
#    graphicFrame (shape 3) from {0B97914B-731A-49B9-966B-491B9EDB8277} to
#    {7DD10862-8B04-4341-A452-6986F8DDAEED}.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{7DD10862-8B04-4341-A452-6986F8DDAEED}")

# ---------------------------------------------------------------------------
# 2) Re-colour the deck's applied theme from the "Integral" palette to the
#    stock "Office Theme" palette (the design/theme swap recorded in the
#    source diff). The ThemeColorScheme index order is fixed:
#    1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6 11=hlink 12=folHlink
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $themeColorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
